$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

# Row 132: H132=3057.4167, I132=3062.0908, J132=3006, K132=9186.2724, L132=9018, M132=-6656.2724, N132=-14078
$ws.Range("H132").Value = 3057.4167
$ws.Range("I132").Value = 3062.0908
$ws.Range("J132").Value = 3006
$ws.Range("K132").Value = 9186.2724
$ws.Range("L132").Value = 9018
$ws.Range("M132").Value = -6656.2724
$ws.Range("N132").Value = -14078

# Row 137: H137=9401.817999999999, I137=6802.5, J137=16333.333, K137=20407.5, L137=48999.999, M137=-17857.5, N137=-54099.999
$ws.Range("H137").Value = 9401.817999999999
$ws.Range("I137").Value = 6802.5
$ws.Range("J137").Value = 16333.333
$ws.Range("K137").Value = 20407.5
$ws.Range("L137").Value = 48999.999
$ws.Range("M137").Value = -17857.5
$ws.Range("N137").Value = -54099.999


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

# Row 132: H132=4442.375, I132=734.36365, K132=2203.09095, M132=326.9090500000002
$ws.Range("H132").Value = 4442.375
$ws.Range("I132").Value = 734.36365
$ws.Range("K132").Value = 2203.09095
$ws.Range("M132").Value = 326.9090500000002


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")

# Row 134: H134=6926.5557, I134=3191.2856, K134=9573.856800000001, M134=-7038.856800000001
$ws.Range("H134").Value = 6926.5557
$ws.Range("I134").Value = 3191.2856
$ws.Range("K134").Value = 9573.856800000001
$ws.Range("M134").Value = -7038.856800000001


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

# Row 2: H2=531.3333, J2=795, L2=795, N2=-1021
$ws.Range("H2").Value = 531.3333
$ws.Range("J2").Value = 795
$ws.Range("L2").Value = 795
$ws.Range("N2").Value = -1021

# Row 5: H5=144.75, I5=93.333336, K5=93.333336, M5=18.666664
$ws.Range("H5").Value = 144.75
$ws.Range("I5").Value = 93.333336
$ws.Range("K5").Value = 93.333336
$ws.Range("M5").Value = 18.666664

# Row 8: H8=1047, J8=100, L8=100, N8=-380
$ws.Range("H8").Value = 1047
$ws.Range("J8").Value = 100
$ws.Range("L8").Value = 100
$ws.Range("N8").Value = -380

# Row 10: H10=577.5, I10=236.25, J10=918.75, K10=236.25, L10=918.75, M10=-97.25, N10=-1196.75
$ws.Range("H10").Value = 577.5
$ws.Range("I10").Value = 236.25
$ws.Range("J10").Value = 918.75
$ws.Range("K10").Value = 236.25
$ws.Range("L10").Value = 918.75
$ws.Range("M10").Value = -97.25
$ws.Range("N10").Value = -1196.75

# Row 11: H11=75, I11=0, J11=75, K11=0, L11=75, N11=-355
$ws.Range("H11").Value = 75
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 75
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 75
$ws.Range("N11").Value = -355
# Row 11: clear M11 (was 40)
$ws.Range("M11").ClearContents()

# Row 12: H12=527.5, J12=527.5, L12=527.5, N12=-867.5
$ws.Range("H12").Value = 527.5
$ws.Range("J12").Value = 527.5
$ws.Range("L12").Value = 527.5
$ws.Range("N12").Value = -867.5

# Row 13: H13=805.8, J13=805.8, L13=805.8, N13=-1083.8
$ws.Range("H13").Value = 805.8
$ws.Range("J13").Value = 805.8
$ws.Range("L13").Value = 805.8
$ws.Range("N13").Value = -1083.8

# Row 14: H14=621, I14=700, J14=601.25, K14=700, L14=601.25, M14=-530, N14=-941.25
$ws.Range("H14").Value = 621
$ws.Range("I14").Value = 700
$ws.Range("J14").Value = 601.25
$ws.Range("K14").Value = 700
$ws.Range("L14").Value = 601.25
$ws.Range("M14").Value = -530
$ws.Range("N14").Value = -941.25

# Row 15: H15=912.8570999999999, I15=1233.3334, J15=672.5, K15=1233.3334, L15=672.5, M15=-1063.3334, N15=-1012.5
$ws.Range("H15").Value = 912.8570999999999
$ws.Range("I15").Value = 1233.3334
$ws.Range("J15").Value = 672.5
$ws.Range("K15").Value = 1233.3334
$ws.Range("L15").Value = 672.5
$ws.Range("M15").Value = -1063.3334
$ws.Range("N15").Value = -1012.5

# Row 21: H21=550, J21=1000, L21=1000, N21=-1470
$ws.Range("H21").Value = 550
$ws.Range("J21").Value = 1000
$ws.Range("L21").Value = 1000
$ws.Range("N21").Value = -1470

# Row 33: H33=2133, I33=1699.5, K33=1699.5, M33=-1320.5
$ws.Range("H33").Value = 2133
$ws.Range("I33").Value = 1699.5
$ws.Range("K33").Value = 1699.5
$ws.Range("M33").Value = -1320.5

# Row 35: H35=1313.2, I35=1313.2, K35=1313.2, M35=-1019.2
$ws.Range("H35").Value = 1313.2
$ws.Range("I35").Value = 1313.2
$ws.Range("K35").Value = 1313.2
$ws.Range("M35").Value = -1019.2

# Row 36: H36=4333, I36=3999.5, J36=5000, K36=3999.5, L36=5000, M36=-3611.5, N36=-5776
$ws.Range("H36").Value = 4333
$ws.Range("I36").Value = 3999.5
$ws.Range("J36").Value = 5000
$ws.Range("K36").Value = 3999.5
$ws.Range("L36").Value = 5000
$ws.Range("M36").Value = -3611.5
$ws.Range("N36").Value = -5776

# Row 37: H37=7525.5, I37=51, J37=15000, K37=51, L37=15000, M37=56, N37=-15214
$ws.Range("H37").Value = 7525.5
$ws.Range("I37").Value = 51
$ws.Range("J37").Value = 15000
$ws.Range("K37").Value = 51
$ws.Range("L37").Value = 15000
$ws.Range("M37").Value = 56
$ws.Range("N37").Value = -15214

# Row 38: H38=1038, I38=1038, K38=1038, M38=-661
$ws.Range("H38").Value = 1038
$ws.Range("I38").Value = 1038
$ws.Range("K38").Value = 1038
$ws.Range("M38").Value = -661

# Row 40: H40=4333, I40=3999.5, J40=5000, K40=3999.5, L40=5000, M40=-3839.5, N40=-5320
$ws.Range("H40").Value = 4333
$ws.Range("I40").Value = 3999.5
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 3999.5
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -3839.5
$ws.Range("N40").Value = -5320

# Row 46: H46=1038, I46=1038, K46=1038, M46=-827
$ws.Range("H46").Value = 1038
$ws.Range("I46").Value = 1038
$ws.Range("K46").Value = 1038
$ws.Range("M46").Value = -827

# Row 56: H56=7000, I56=7000, K56=7000, M56=-6155
$ws.Range("H56").Value = 7000
$ws.Range("I56").Value = 7000
$ws.Range("K56").Value = 7000
$ws.Range("M56").Value = -6155

# Row 58: H58=6753, I58=1506, K58=1506, M58=-1303
$ws.Range("H58").Value = 6753
$ws.Range("I58").Value = 1506
$ws.Range("K58").Value = 1506
$ws.Range("M58").Value = -1303

# Row 86: H86=4666.3335, I86=3999.5, K86=3999.5, M86=-2876.5
$ws.Range("H86").Value = 4666.3335
$ws.Range("I86").Value = 3999.5
$ws.Range("K86").Value = 3999.5
$ws.Range("M86").Value = -2876.5

# Row 89: H89=4666.3335, I89=3999.5, K89=19997.5, M89=-14381.5
$ws.Range("H89").Value = 4666.3335
$ws.Range("I89").Value = 3999.5
$ws.Range("K89").Value = 19997.5
$ws.Range("M89").Value = -14381.5

# Row 136: H136=6753, I136=1506, K136=4518, M136=-1968
$ws.Range("H136").Value = 6753
$ws.Range("I136").Value = 1506
$ws.Range("K136").Value = 4518
$ws.Range("M136").Value = -1968


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

# Row 126: H126=7757, I126=0, J126=7757, K126=0, L126=23271, N126=-28211
$ws.Range("H126").Value = 7757
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 7757
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 23271
$ws.Range("N126").Value = -28211
# Row 126: clear M126 (was -27530)
$ws.Range("M126").ClearContents()


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

# Row 22: H22=2643.2856, J22=3380.4, L22=3380.4, N22=-3970.4
$ws.Range("H22").Value = 2643.2856
$ws.Range("J22").Value = 3380.4
$ws.Range("L22").Value = 3380.4
$ws.Range("N22").Value = -3970.4

# Row 27: H27=2643.2856, J27=3380.4, L27=3380.4, N27=-3594.4
$ws.Range("H27").Value = 2643.2856
$ws.Range("J27").Value = 3380.4
$ws.Range("L27").Value = 3380.4
$ws.Range("N27").Value = -3594.4

# Row 46: H46=6170.5884, J46=6062.5, L46=6062.5, N46=-6438.5
$ws.Range("H46").Value = 6170.5884
$ws.Range("J46").Value = 6062.5
$ws.Range("L46").Value = 6062.5
$ws.Range("N46").Value = -6438.5

# Row 55: H55=979.9, I55=949.8570999999999, J55=1050, K55=949.8570999999999, L55=1050, M55=-776.8570999999999, N55=-1396
$ws.Range("H55").Value = 979.9
$ws.Range("I55").Value = 949.8570999999999
$ws.Range("J55").Value = 1050
$ws.Range("K55").Value = 949.8570999999999
$ws.Range("L55").Value = 1050
$ws.Range("M55").Value = -776.8570999999999
$ws.Range("N55").Value = -1396

# Row 132: H132=17878.223, J132=20780, L132=62340, N132=-67400
$ws.Range("H132").Value = 17878.223
$ws.Range("J132").Value = 20780
$ws.Range("L132").Value = 62340
$ws.Range("N132").Value = -67400

